# "compara en casa terminado" - replace the 123seguro lead URLs with the
# comparaencasa.com quote/result URLs, add a "Fecha" column, and rename
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet name: "urls" -> "Sheet1" ---
$ws.Name = "Sheet1"

# --- drop the stray formatted-but-empty D12 cell / row left over at the
#     bottom of the sheet (also restores dimension to the real used range)
$ws.Rows.Item(12).Delete()

# --- new "Fecha" header in C1, matching the existing header style (B1) ---
$ws.Cells.Item(1, 3).Value = "Fecha"
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)   # xlPasteFormats

# --- new URLs for column B (ids in column A are unchanged) ---
$ws.Cells.Item(2, 2).Value = "https://www.comparaencasa.com/seguros-de-auto/cotizador/?car_brand=28&year=2020"
$ws.Cells.Item(3, 2).Value = "https://www.comparaencasa.com/seguros-de-auto/resultados/ford-en-focus-l-16-1-6-4-p-s-capital-federal/c7a1294e28bc43ad8bfda0387f413c33/"
$ws.Cells.Item(4, 2).Value = "https://www.comparaencasa.com/seguros-de-auto/resultados/toyota-en-corolla-2-0-xei-l-20-cvt-capital-federal/fa2675522f8c4017a4993887f7f7d4a1/"
$ws.Cells.Item(5, 2).Value = "https://www.comparaencasa.com/seguros-de-auto/cotizador/?car_brand=153&year=2019"
$ws.Cells.Item(6, 2).Value = "https://www.comparaencasa.com/seguros-de-auto/cotizador/?car_brand=28&year=2018"

# --- new "Fecha" values for column C ---
$ws.Cells.Item(2, 3).Value = "15-06-2021"
$ws.Cells.Item(3, 3).Value = "15-06-2021"
$ws.Cells.Item(4, 3).Value = "15-06-2021"
$ws.Cells.Item(5, 3).Value = "15-06-2021"
$ws.Cells.Item(6, 3).Value = "15-06-2021"

# --- rebuild hyperlinks for the new URLs (old ones pointed at 123seguro) ---
# stash the existing hyperlink cell style off-sheet first: Hyperlinks.Add()
# reformats the target cell's font, so we restore the original look after.
$ws.Cells.Item(3, 2).Copy()
$ws.Cells.Item(100, 100).PasteSpecial(-4122)

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "https://www.comparaencasa.com/seguros-de-auto/cotizador/?car_brand=28&year=2020")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 2), "https://www.comparaencasa.com/seguros-de-auto/resultados/ford-en-focus-l-16-1-6-4-p-s-capital-federal/c7a1294e28bc43ad8bfda0387f413c33/")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 2), "https://www.comparaencasa.com/seguros-de-auto/resultados/toyota-en-corolla-2-0-xei-l-20-cvt-capital-federal/fa2675522f8c4017a4993887f7f7d4a1/")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 2), "https://www.comparaencasa.com/seguros-de-auto/cotizador/?car_brand=153&year=2019")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 2), "https://www.comparaencasa.com/seguros-de-auto/cotizador/?car_brand=28&year=2018")

$ws.Cells.Item(100, 100).Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)
$ws.Cells.Item(100, 100).Clear()

# --- leave the cursor on A1 (the stray D12 selection from the source file
#     no longer points at real data now that row 12 is gone) ---
$ws.Range("A1").Select()

